$wb = $excel.ActiveWorkbook

# --- Workbook-level view settings ---
$excel.ActiveWindow.Width = 14120

# --- metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Range("D6").Value = 35
$wsMeta.Range("D7").Value = 0
$wsMeta.Range("G7").Value = 81
$wsMeta.Range("H7").Value = 66
$wsMeta.Range("D8").Select()

# --- depths sheet updates ---
$wsDepths = $wb.Worksheets.Item("depths")

$data = @{
    2 = @{ J=127.715; K=127.725; L=127.68; M=127.694 }
    3 = @{ J=127.725; K=127.736; L=127.694; M=127.709 }
    4 = @{ J=127.736; K=127.747; L=127.709; M=127.724 }
    5 = @{ J=127.747; K=127.759; L=127.724; M=127.739 }
    6 = @{ J=127.759; K=127.77; L=127.739; M=127.754 }
    7 = @{ J=127.77; K=127.781; L=127.754; M=127.769 }
    8 = @{ J=127.781; K=127.792; L=127.769; M=127.784 }
    9 = @{ J=127.792; K=127.804; L=127.784; M=127.799 }
    10 = @{ J=127.804; K=127.815; L=127.799; M=127.814 }
    11 = @{ J=127.815; K=127.826; L=127.814; M=127.829 }
    12 = @{ J=127.826; K=127.837; L=127.829; M=127.844 }
    13 = @{ J=127.837; K=127.848; L=127.844; M=127.86 }
    14 = @{ J=127.848; K=127.86; L=127.86; M=127.875 }
    15 = @{ J=127.86; K=127.871; L=127.875; M=127.89 }
    16 = @{ J=127.871; K=127.882; L=127.89; M=127.905 }
    17 = @{ J=127.882; K=127.893; L=127.905; M=127.92 }
    18 = @{ J=127.893; K=127.904; L=127.92; M=127.935 }
    19 = @{ J=127.904; K=127.916; L=127.935; M=127.95 }
    20 = @{ J=127.916; K=127.927; L=127.95; M=127.965 }
    21 = @{ J=127.927; K=127.938; L=127.965; M=127.98 }
    22 = @{ J=127.938; K=127.949; L=127.98; M=127.995 }
    23 = @{ J=127.949; K=127.961; L=127.995; M=128.01 }
    24 = @{ J=127.961; K=127.972; L=128.01; M=128.025 }
    25 = @{ J=127.972; K=127.983; L=128.025; M=128.04 }
    26 = @{ J=127.983; K=127.994; L=128.04; M=128.055 }
    27 = @{ J=127.994; K=128.005; L=128.055; M=128.07 }
    28 = @{ J=128.005; K=128.017; L=128.07; M=128.085 }
    29 = @{ J=128.017; K=128.028; L=128.085; M=128.1 }
    30 = @{ J=128.028; K=128.039; L=128.1; M=128.115 }
    31 = @{ J=128.039; K=128.05; L=128.115; M=128.13 }
    32 = @{ J=128.05; K=128.061; L=128.13; M=128.145 }
    33 = @{ J=128.061; K=128.073; L=128.145; M=128.16 }
    34 = @{ J=128.073; K=128.084; L=128.16; M=128.175 }
    35 = @{ J=128.084; K=128.095; L=128.175; M=128.191 }
    36 = @{ J=128.095; K=128.106; L=128.191; M=128.206 }
    37 = @{ J=128.106; K=128.118; L=128.206; M=128.221 }
    38 = @{ J=128.118; K=128.129; L=128.221; M=128.236 }
    39 = @{ J=128.129; K=128.14; L=128.236; M=128.251 }
    40 = @{ J=128.14; K=128.151; L=128.251; M=128.266 }
    41 = @{ J=128.151; K=128.162; L=128.266; M=128.281 }
    42 = @{ J=128.162; K=128.174; L=128.281; M=128.296 }
    43 = @{ J=128.174; K=128.185; L=128.296; M=128.311 }
    44 = @{ J=128.185; K=128.196; L=128.311; M=128.326 }
    45 = @{ J=128.196; K=128.207; L=128.326; M=128.341 }
    46 = @{ J=128.207; K=128.218; L=128.341; M=128.364 }
    47 = @{ J=128.218; K=128.23 }
    48 = @{ J=128.23; K=128.241 }
    49 = @{ J=128.241; K=128.252 }
    50 = @{ J=128.252; K=128.263 }
    51 = @{ J=128.263; K=128.275 }
    52 = @{ J=128.275; K=128.286 }
    53 = @{ J=128.286; K=128.297 }
    54 = @{ J=128.297; K=128.308 }
    55 = @{ J=128.308; K=128.319 }
    56 = @{ J=128.319; K=128.331 }
    57 = @{ J=128.331; K=128.342 }
    58 = @{ J=128.342; K=128.36 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $wsDepths.Range($addr).Value = $cols[$col]
    }
}

# Re-apply header/index-column formatting (bold "Aptos Narrow" 12pt, centered/top, boxed)
$hdrRange = $wsDepths.Range("B1:U1")
$idxRange = $wsDepths.Range("A2:A61")
$fmtRange = $wsDepths.Range("B1:U1,A2:A61")

$fmtRange.Borders.Color = 0
$fmtRange.Font.Name = "Aptos Narrow"
$fmtRange.Font.Size = 12
$fmtRange.Font.Bold = $true
$fmtRange.HorizontalAlignment = -4108
$fmtRange.VerticalAlignment = -4160

# --- sample_thicknesses sheet view selection ---
$wsThick = $wb.Worksheets.Item("sample_thicknesses")
$wsThick.Range("C49").Select()

Write-Host "edit complete"
